$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestResults")

# Data rows: Test Case ID, Test scenario, Status
$data = @(
    @("AL_001",  "@API Admin Login Success with Mandatory Fields", "passed"),
    @("QS_067",  "@API Validation of ISAWE-CASE question successfull message.", "passed"),
    @("QS_068",  "@API Validation of edit ISAWE-CASE question successfull message.", "passed"),
    @("QS_069",  "@API Validation of ISAWE-CASE question Approved message.", "passed"),
    @("QS_070",  "@API Validation of ISAWE-CASE question checkout message.", "passed"),
    @("QS_069A", "@API Validation of ISAWE-CASE question Approved message (again).", "passed"),
    @("QS_071",  "@API ISAWE-CASE endpoint validation", "passed"),
    @("QS_072",  "@API ISAWE-CASE Method validation-  incorrect HTTP method", "passed"),
    @("QS_073",  "@API ISAWE-CASE Header field validation - invalid", "passed"),
    @("QS_074",  "@API Validation of empty title field for ISAWE-CASE", "failed")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
